$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set up column D formatting first, mirroring column C ---
# D1 gets the centered header style used by A1/B1/C1.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
# D2:D54 get the plain centered-text style used by the rest of the table.
$ws.Range("C2:C54").Copy()
$ws.Range("D2:D54").PasteSpecial(-4122)

# Column D width (closest value the engine's column-width grid allows to 27.08984375).
$ws.Columns("D").ColumnWidth = 26.4

# --- New column D values: header (row 1) plus 53 data rows (rows 2-54). ---
# Rows 55-67 intentionally keep only columns A/C, matching the source diff.
$dValues = @(
    "vinhquang13531@gmail.com",
    "492401487499641",
    "264846833886893",
    "246933925685884",
    "275712962530384",
    "179112262432141",
    "264292347041953",
    "781074362021860",
    "319718041464688",
    "108757516583030",
    "1911119432502301",
    "1471037103142820",
    "498723016852762",
    "451122651649814",
    "1455330608067225",
    "459470627407885",
    "171076470160303",
    "154362801673336",
    "703894706306308",
    "1826787617399125",
    "1209295105801230",
    "126348794234977",
    "1761231377490564",
    "145971692474202",
    "205882209595847",
    "761161927309184",
    "423158777729623",
    "774067269282955",
    "1549979828593724",
    "415793858487646",
    "386412421458390",
    "762027117282800",
    "616734738469807",
    "1544909859170915",
    "332252896899497",
    "1856620574550792",
    "1665045370384962",
    "159617931477214",
    "2116731221905896",
    "784012071764293",
    "206240536583250",
    "1384506521580724",
    "494023894278277",
    "737452329661021",
    "1720869341538398",
    "322008967984765",
    "277314789020591",
    "1798895427089724",
    "1208088295924736",
    "316603388502163",
    "249978345123454",
    "436038656468443",
    "173960829412142",
    "1362279090456478"
)

for ($i = 0; $i -lt $dValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 4).Value = $dValues[$i]
}

# D1 is a mailto: hyperlink, just like the A1/B1/C1 header cells.
$ws.Hyperlinks.Add($ws.Cells.Item(1, 4), "mailto:vinhquang13531@gmail.com")

# Hyperlinks.Add re-stamps its own style onto the cell, so restore the header
# style/value to match C1 after adding the link.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Cells.Item(1, 4).Value = $dValues[0]

# Active cell shown when the workbook is opened.
$ws.Range("G7").Select()
